# Hortaliza, Macroferia Regional de Talca - Zapallo
# A new weekly price observation (Camote, "1a (guarda)") is inserted as a
# new data row right after the header block of existing rows, at sheet
# row 208. This pushes the previously-existing rows 208-244 down to
# rows 209-245 (dimension grows from A1:R244 to A1:R245), exactly as a
# native Excel "insert row" does - all other rows/columns are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 208, shifting rows 208:244 down to 209:245.
$ws.Rows.Item(208).Insert()

# Populate the newly inserted row 208 with the new observation.
$ws.Cells.Item(208, 1).Value  = 5
$ws.Cells.Item(208, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(208, 3).Value  = "Maule"
$ws.Cells.Item(208, 4).Value  = 44694
$ws.Cells.Item(208, 5).Value  = 7
$ws.Cells.Item(208, 6).Value  = 100112045
$ws.Cells.Item(208, 7).Value  = "Zapallo"
$ws.Cells.Item(208, 8).Value  = "Camote"
$ws.Cells.Item(208, 9).Value  = "1a (guarda)"
$ws.Cells.Item(208, 10).Value = 900
$ws.Cells.Item(208, 11).Value = 350
$ws.Cells.Item(208, 12).Value = 350
$ws.Cells.Item(208, 13).Value = 350
$ws.Cells.Item(208, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(208, 15).Value = "Región del Maule"
$ws.Cells.Item(208, 16).Value = 350
$ws.Cells.Item(208, 17).Value = 1
$ws.Cells.Item(208, 18).Value = "Hortaliza"
